# Update the "mail" header to "email" on the user sheet, and move the
# active selection from G6 to E5 (matching the author's edit in Excel).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D1 currently holds "mail" -> rename the column header to "email".
$ws.Range("D1").Value = "email"

# Update the current selection to E5, as captured in the saved view state.
$ws.Range("E5").Select()
